$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Shared text used across sheets
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$targetFile1 = "2f80e1f8-af23-4406-a0df-9a04aefade6a.md"
$targetFile2 = "4807d2ba-722e-4846-8e6e-0d411ebd2e7f.md"
$url1 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/39010130ca1cff18aac03a001750dfa93863bcc4/e2e/2f80e1f8-af23-4406-a0df-9a04aefade6a.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/39010130ca1cff18aac03a001750dfa93863bcc4/e2e/4807d2ba-722e-4846-8e6e-0d411ebd2e7f.md"

# ---------------------------------------------------------------------------
# Overview sheet: status text changed (drives column width change via autofit)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.2
$wsOverview.Columns.Item(6).ColumnWidth = 29.2

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus
$wsZh.Columns.Item(3).ColumnWidth = 29.2

# Latest Target File / Latest Handback File
$wsZh.Range("I2").Value = $targetFile1
$wsZh.Range("J2").Value = "2f80e1f8-af23-4406-a0df-9a04aefade6a.8e8d1c9f5056664513c1a7a9692ea42b7657589b.zh-cn.xlf"
$wsZh.Range("I3").Value = $targetFile2
$wsZh.Range("J3").Value = "4807d2ba-722e-4846-8e6e-0d411ebd2e7f.7b00e1d0f351fef60032176e58d4325edc14803b.zh-cn.xlf"
$wsZh.Columns.Item(9).ColumnWidth = 39.2
$wsZh.Columns.Item(10).ColumnWidth = 39.2

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $url1, [Type]::Missing, [Type]::Missing, $targetFile1)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $url2, [Type]::Missing, [Type]::Missing, $targetFile2)

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus
$wsDe.Columns.Item(3).ColumnWidth = 29.2

$wsDe.Range("I2").Value = $targetFile1
$wsDe.Range("J2").Value = "2f80e1f8-af23-4406-a0df-9a04aefade6a.8e8d1c9f5056664513c1a7a9692ea42b7657589b.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-14 03:06:49"
$wsDe.Range("I3").Value = $targetFile2
$wsDe.Range("J3").Value = "4807d2ba-722e-4846-8e6e-0d411ebd2e7f.7b00e1d0f351fef60032176e58d4325edc14803b.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-14 03:06:49"
$wsDe.Columns.Item(9).ColumnWidth = 39.2
$wsDe.Columns.Item(10).ColumnWidth = 39.2

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $url1, [Type]::Missing, [Type]::Missing, $targetFile1)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $url2, [Type]::Missing, [Type]::Missing, $targetFile2)
